$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: MuSCs / Bdnf / Sort1 -> ECs
$ws.Cells.Item(2, 1).Value = "MuSCs"
$ws.Cells.Item(2, 2).Value = "Bdnf"
$ws.Cells.Item(2, 3).Value = "Sort1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.9253576666666667
$ws.Cells.Item(2, 8).Value = 2.776073
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.4652636666666667
$ws.Cells.Item(2, 14).Value = 1.395791
$ws.Cells.Item(2, 15).Value = 0.02604271297411062
$ws.Cells.Item(2, 16).Value = 0.02604271297411062
$ws.Cells.Item(2, 17).Value = 0.4305353009714445
$ws.Cells.Item(2, 18).Value = 3.874817708743
$ws.Cells.Item(2, 19).Value = 0.02604271297411062
$ws.Cells.Item(2, 20).Value = 0.02604271297411062

# Row 3: MuSCs / Bdnf / Sort1 -> FAPs
$ws.Cells.Item(3, 1).Value = "MuSCs"
$ws.Cells.Item(3, 2).Value = "Bdnf"
$ws.Cells.Item(3, 3).Value = "Sort1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.9253576666666667
$ws.Cells.Item(3, 8).Value = 2.776073
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.781382333333333
$ws.Cells.Item(3, 14).Value = 5.344147
$ws.Cells.Item(3, 15).Value = 0.09971126509087273
$ws.Cells.Item(3, 16).Value = 0.09971126509087272
$ws.Cells.Item(3, 17).Value = 1.648415799414556
$ws.Cells.Item(3, 18).Value = 14.835742194731
$ws.Cells.Item(3, 19).Value = 0.09971126509087273
$ws.Cells.Item(3, 20).Value = 0.09971126509087272

# Row 4: MuSCs / Bdnf / Sort1 -> MuSCs
$ws.Cells.Item(4, 1).Value = "MuSCs"
$ws.Cells.Item(4, 2).Value = "Bdnf"
$ws.Cells.Item(4, 3).Value = "Sort1"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.9253576666666667
$ws.Cells.Item(4, 8).Value = 2.776073
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 15.618761
$ws.Cells.Item(4, 14).Value = 46.856283
$ws.Cells.Item(4, 15).Value = 0.8742460219350168
$ws.Cells.Item(4, 16).Value = 0.8742460219350167
$ws.Cells.Item(4, 17).Value = 14.45294023518433
$ws.Cells.Item(4, 18).Value = 130.076462116659
$ws.Cells.Item(4, 19).Value = 0.8742460219350168
$ws.Cells.Item(4, 20).Value = 0.8742460219350167

# Remove the now-duplicate trailing rows (5:7) so the sheet shrinks to A1:T4
$ws.Rows("5:7").Delete() | Out-Null
